$d = $word.ActiveDocument

# Locate the paragraph that contains the "Invalid block" error run
# (the second paragraph of the extraSpaceInEndTemplate sample). We search
# by content instead of a hard-coded index so the script is resilient to
# minor structural differences.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text.Contains("Invalid block") -and $text.Contains("[ENDTEMPLATE]")) {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # Rebuild the paragraph keeping only the bold/red error run, dropping
    # the leading empty run (<w:t/>) and its <w:pPr><w:rPr><w:lang .../>
    # that used to precede it (this was the extra, unwanted run left over
    # from sub-template generation - see #348).
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r>' +
             '<w:rPr>' +
               '<w:b w:val="true"/>' +
               '<w:color w:val="FF0000"/>' +
             '</w:rPr>' +
             '<w:t>Invalid block: Unexpected tag EOF missing [ENDTEMPLATE]</w:t>' +
           '</w:r>' +
         '</w:p>'

    $r.InsertXML($xml)
}
